$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Commercial revenue indicator): clarify units as "millions" of 2024 USD
$ws.Range("A4").Value = "Commercial revenue (millions 2024 USD)"

# Row 6 (Bottom temperature): replace implications text - updated cold pool language + source link
$ws.Range("C6").Value = "Inshore temperature thresholds (around 14°C) initiate migration of squid from offshore overwintering habitats. Longfin squid seasonal distribution and growth rates are likely temperature dependent, avoiding water <8°C. Stronger and/or more persistent Mid-Atlantic Cold Pool conditions (not shown) may limit habitat availability (https://noaa-edab.github.io/catalog/cold_pool.html)."

# Row 5 (Western Gulf Stream Index): append source link to the implications text
$ws.Range("C5").Value = "Since the mid-1990s, north and westward shifts in the Gulf Stream have resulted in an increase in warm core rings and deep water, high salinity heat waves. The position of the Gulf Stream influences seasonal temperature and water mass mixing dynamics that affect longfin squid habitat suitability, temperature-dependent growth, and prey availability (https://noaa-edab.github.io/catalog/gsi.html). "

# Update the active selection saved in the sheet view
$ws.Range("D11").Select()
